# December-2020-Collection.xlsx update
# "data updated till 15Dec 9AM"
# Adds a new day's column (S = 14-Dec-2020) of collection figures across
# many rows, updates the per-row totals (column E, formula-driven, will
# recalc automatically) and adds three new comments from "Vijay" on some
# of the new cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Purple highlight color used for some of the "S" cells (matches the
# existing fill already used elsewhere in the sheet, e.g. N3/I10/etc.)
$purple = 9660795

# row -> (value, useHighlight)
$rows = @(
    @(3,  5000, $true),
    @(5,  1000, $true),
    @(10, 2000, $false),
    @(13, 3000, $false),
    @(14, 2000, $true),
    @(16, 2000, $false),
    @(21, 1000, $false),
    @(25, 3000, $false),
    @(30, 2000, $false),
    @(33, 3000, $false),
    @(35, 2000, $true),
    @(43, 1500, $false),
    @(54, 3000, $true),
    @(57, 1000, $true),
    @(59, 4500, $true),
    @(61, 2000, $false),
    @(63, 2000, $false),
    @(65, 5000, $false),
    @(66, 2000, $false),
    @(68, 3000, $true),
    @(70, 3000, $false),
    @(71, 3000, $true),
    @(79, 2000, $true),
    @(80, 5000, $false),
    @(81, 2000, $false),
    @(86, 3000, $false),
    @(89, 5000, $false),
    @(90, 1000, $false),
    @(91, 1000, $true)
)

foreach ($row in $rows) {
    $r = $row[0]
    $val = $row[1]
    $hl = $row[2]
    $cell = $ws.Range("S" + $r)
    $cell.Value = $val
    if ($hl) {
        $cell.Interior.Color = $purple
    }
}

# New comments left by Vijay on three of the new cells
function Add-VijayComment($addr, $line2) {
    $cmt = $ws.Range($addr).AddComment()
    $full = "Vijay:" + [char]10 + $line2
    $cmt.Text($full)
    $tf = $cmt.Shape.TextFrame
    $bold = $tf.Characters(1, 6)
    $bold.Font.Bold = $true
    $bold.Font.Name = "Times New Roman"
    $bold.Font.Size = 9
    $rest = $tf.Characters(7, 200)
    $rest.Font.Bold = $false
    $rest.Font.Name = "Times New Roman"
    $rest.Font.Size = 9
}

Add-VijayComment "S54" ("2500-Digital" + [char]10 + "500-Cash")
Add-VijayComment "S68" ("2000-Cash" + [char]10 + "1000-Digital")
Add-VijayComment "S79" ("1000-Cash" + [char]10 + "1000-Digital")

# Move the frozen-pane cursor to where the user was last working
$ws.Range("S64").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 13
